$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 54: Arcane Arts for Dummies | Book of Mythril
$ws.Range("H54").Value = 50000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
# ALC row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Range("H106").Value = 2920
$ws.Range("I106").Value = 2760
$ws.Range("J106").Value = 3400
$ws.Range("K106").Value = 2760
$ws.Range("L106").Value = 3400
$ws.Range("M106").Value = -2129
$ws.Range("N106").Value = -4662
# ALC row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 364135.97
$ws.Range("I116").Value = 716071.4399999999
$ws.Range("K116").Value = 716071.4399999999
$ws.Range("M116").Value = -712629.4399999999
# ALC row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 2495.4285
$ws.Range("I125").Value = 2844
$ws.Range("J125").Value = 2234
$ws.Range("K125").Value = 25596
$ws.Range("L125").Value = 20106
$ws.Range("M125").Value = -23136
$ws.Range("N125").Value = -25026
# ALC row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 948.53845
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 959.2941
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 2877.8823
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -12877.8823

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 1162
$ws.Range("I2").Value = 1211.3
$ws.Range("J2").Value = 997.6667
$ws.Range("K2").Value = 1211.3
$ws.Range("L2").Value = 997.6667
$ws.Range("M2").Value = -1098.3
$ws.Range("N2").Value = -1223.6667
# ARM row 25: Still Crazy After All These Years | Initiate's Alembic
$ws.Range("H25").Value = 3903.2
$ws.Range("I25").Value = 2129
$ws.Range("K25").Value = 2129
$ws.Range("M25").Value = -1727
# ARM row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 2698.6191
$ws.Range("I74").Value = 2112.3333
$ws.Range("J74").Value = 4164.3335
$ws.Range("K74").Value = 2112.3333
$ws.Range("L74").Value = 4164.3335
$ws.Range("M74").Value = -1238.3333
$ws.Range("N74").Value = -5912.3335
# ARM row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 2698.6191
$ws.Range("I77").Value = 2112.3333
$ws.Range("J77").Value = 4164.3335
$ws.Range("K77").Value = 10561.6665
$ws.Range("L77").Value = 20821.6675
$ws.Range("M77").Value = -6193.666499999999
$ws.Range("N77").Value = -29557.6675
# ARM row 95: Shielded Life | High Steel Scutum
$ws.Range("H95").Value = 49900
$ws.Range("J95").Value = 49900
$ws.Range("L95").Value = 49900
$ws.Range("N95").Value = -55392
# ARM row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 1162
$ws.Range("I116").Value = 1211.3
$ws.Range("J116").Value = 997.6667
$ws.Range("K116").Value = 1211.3
$ws.Range("L116").Value = 997.6667
$ws.Range("M116").Value = 1082.7
$ws.Range("N116").Value = -5585.6667
# ARM row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2328.1052
$ws.Range("I132").Value = 1089
$ws.Range("J132").Value = 4711
$ws.Range("K132").Value = 3267
$ws.Range("L132").Value = 14133
$ws.Range("M132").Value = -737
$ws.Range("N132").Value = -19193

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 1162
$ws.Range("I3").Value = 1211.3
$ws.Range("J3").Value = 997.6667
$ws.Range("K3").Value = 1211.3
$ws.Range("L3").Value = 997.6667
$ws.Range("M3").Value = -1097.3
$ws.Range("N3").Value = -1225.6667
# BSM row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 1773.409
$ws.Range("I94").Value = 1821.8422
$ws.Range("J94").Value = 1466.6666
$ws.Range("K94").Value = 1821.8422
$ws.Range("L94").Value = 1466.6666
$ws.Range("M94").Value = -1370.8422
$ws.Range("N94").Value = -2368.6666
# BSM row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2580
$ws.Range("I105").Value = 2800
$ws.Range("J105").Value = 2066.6667
$ws.Range("K105").Value = 2800
$ws.Range("L105").Value = 2066.6667
$ws.Range("M105").Value = -1053
$ws.Range("N105").Value = -5560.6667
# BSM row 128: Mangalomania | Manganese Ingot
$ws.Range("H128").Value = 1840
$ws.Range("I128").Value = 1840
$ws.Range("K128").Value = 5520
$ws.Range("M128").Value = -3030

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 62.214287
$ws.Range("I7").Value = 62.583332
$ws.Range("K7").Value = 62.583332
$ws.Range("M7").Value = 50.416668
# CRP row 69: Landing the Big One | Cedar Fishing Rod
$ws.Range("H69").Value = 24242.857
$ws.Range("I69").Value = 14940
$ws.Range("K69").Value = 14940
$ws.Range("M69").Value = -14191
# CRP row 72: Fishing for Profits (L) | Cedar Fishing Rod
$ws.Range("H72").Value = 24242.857
$ws.Range("I72").Value = 14940
$ws.Range("K72").Value = 44820
$ws.Range("M72").Value = -41076
# CRP row 93: Reeling for Rods | Muudhorn Fishing Rod
$ws.Range("H93").Value = 8204.5
$ws.Range("I93").Value = 8204.5
$ws.Range("K93").Value = 8204.5
$ws.Range("M93").Value = -6332.5
# CRP row 103: Spare a Rod and Spoil the Fishers | Gazelle Horn Fishing Rod
$ws.Range("H103").Value = 19805
$ws.Range("I103").Value = 9006.25
$ws.Range("J103").Value = 63000
$ws.Range("K103").Value = 9006.25
$ws.Range("L103").Value = 63000
$ws.Range("M103").Value = -7834.25
$ws.Range("N103").Value = -65344
# CRP row 112: Understaffed | Applewood Cane
$ws.Range("H112").Value = 27425
$ws.Range("J112").Value = 27425
$ws.Range("L112").Value = 27425
$ws.Range("N112").Value = -30379
# CRP row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 5493
$ws.Range("I134").Value = 6041.95
$ws.Range("K134").Value = 18125.85
$ws.Range("M134").Value = -15590.85

$ws = $wb.Worksheets.Item("CUL")
# CUL row 141: Ocean Explosion | Acqua Pazza
$ws.Range("H141").Value = 6087.4165
$ws.Range("I141").Value = 5731.727
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 17195.181
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -12015.181
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
# GSM row 19: Better Four Eyes than None | Brass Spectacles
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
# GSM row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 2685
$ws.Range("I102").Value = 1996.5555
$ws.Range("J102").Value = 5783
$ws.Range("K102").Value = 1996.5555
$ws.Range("L102").Value = 5783
$ws.Range("M102").Value = -374.5554999999999
$ws.Range("N102").Value = -9027
# GSM row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 10101677
$ws.Range("I107").Value = 678.4286
$ws.Range("J107").Value = 27778424
$ws.Range("K107").Value = 678.4286
$ws.Range("L107").Value = 27778424
$ws.Range("M107").Value = 1241.5714
$ws.Range("N107").Value = -27782264
# GSM row 119: Bulking Up | Dwarven Mythril Rapier
$ws.Range("H119").Value = 39766.668
$ws.Range("J119").Value = 39766.668
$ws.Range("L119").Value = 39766.668
$ws.Range("N119").Value = -49442.668

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 2274.9167
$ws.Range("I22").Value = 2475
$ws.Range("J22").Value = 2174.875
$ws.Range("K22").Value = 2475
$ws.Range("L22").Value = 2174.875
$ws.Range("M22").Value = -2180
$ws.Range("N22").Value = -2764.875
# LTW row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 2274.9167
$ws.Range("I27").Value = 2475
$ws.Range("J27").Value = 2174.875
$ws.Range("K27").Value = 2475
$ws.Range("L27").Value = 2174.875
$ws.Range("M27").Value = -2368
$ws.Range("N27").Value = -2388.875
# LTW row 30: Packing a Punch | Goatskin Cesti
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
# LTW row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 2086
$ws.Range("I100").Value = 1932.6364
$ws.Range("J100").Value = 2367.1667
$ws.Range("K100").Value = 1932.6364
$ws.Range("L100").Value = 2367.1667
$ws.Range("M100").Value = -1391.6364
$ws.Range("N100").Value = -3449.1667
# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3872.8367
$ws.Range("I132").Value = 2017.3334
$ws.Range("J132").Value = 4474.6216
$ws.Range("K132").Value = 6052.0002
$ws.Range("L132").Value = 13423.8648
$ws.Range("M132").Value = -3522.0002
$ws.Range("N132").Value = -18483.8648

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 5269.1875
$ws.Range("I122").Value = 3792.8462
$ws.Range("K122").Value = 11378.5386
$ws.Range("M122").Value = -8928.5386
# WVR row 131: A Better Bottom Line | AR-Caean Velvet Bottoms of Scouting
$ws.Range("H131").Value = 51257.777
$ws.Range("J131").Value = 51257.777
$ws.Range("L131").Value = 51257.777
$ws.Range("N131").Value = -61337.777
# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 9264077
$ws.Range("I132").Value = 8475.691999999999
$ws.Range("K132").Value = 25427.076
$ws.Range("M132").Value = -22897.076
